$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.057.04'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '1.709.12'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.00'
$ws.Range('E5').Value = '  -6.14%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '0.4729'
$ws.Range('E7').Value = '  +5.88%  '
$ws.Range('D8').Value = '0.3417'
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('D9').Value = '42.07'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '0.07257'
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('D11').Value = '1.032'
$ws.Range('E11').Value = '  -5.90%  '
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '19.74'
$ws.Range('E13').Value = '  -5.57%  '
$ws.Range('D14').Value = '5.825'
$ws.Range('E14').Value = '  -3.14%  '
$ws.Range('D15').Value = '1.706.92'
$ws.Range('E15').Value = '  -3.57%  '
$ws.Range('D16').Value = '6.819'
$ws.Range('E16').Value = '  -5.68%  '
$ws.Range('D17').Value = '88.82'
$ws.Range('E17').Value = '  -4.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001034'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').Value = '0.06354'
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('E21').Value = '  -3.90%  '
$ws.Range('D22').Value = '5.589'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('D23').Value = '27.085.57'
$ws.Range('E23').Value = '  -2.85%  '
$ws.Range('D24').Value = '10.82'
$ws.Range('E24').Value = '  -3.96%  '
$ws.Range('D25').Value = '2.111'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.50'
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D27').Value = '19.51'
$ws.Range('E27').Value = '  -4.21%  '
$ws.Range('D28').Value = '1.901.49'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').Value = '2.061'
$ws.Range('E29').Value = '  -4.35%  '
$ws.Range('D30').Value = '119.08'
$ws.Range('E30').Value = '  -4.58%  '
$ws.Range('D31').Value = '1.006'
$ws.Range('E31').Value = '  -9.05%  '
$ws.Range('D32').Value = '0.09139'
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').Value = '3.579'
$ws.Range('E33').Value = '  -2.72%  '
$ws.Range('D34').Value = '5.282'
$ws.Range('E34').Value = '  -6.07%  '
$ws.Range('D35').Value = '0.02181'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('D36').Value = '0.05799'
$ws.Range('E36').Value = '  -6.20%  '
$ws.Range('E37').Value = '  -7.12%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.000'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.1979'
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('D40').Value = '4.718'
$ws.Range('E40').Value = '  -4.65%  '
$ws.Range('D41').Value = '1.389'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '0.5859'
$ws.Range('E42').Value = '  -7.04%  '
$ws.Range('D43').Value = '1.096'
$ws.Range('E43').Value = '  -7.40%  '
$ws.Range('D44').Value = '7.441'
$ws.Range('E44').Value = '  -5.30%  '
$ws.Range('E45').Value = '  -5.37%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.557'
$ws.Range('E46').Value = '  -5.12%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5614'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').Value = '117.47'
$ws.Range('E48').Value = '  -4.08%  '
$ws.Range('D49').Value = '1.828'
$ws.Range('E49').Value = '  -6.36%  '
$ws.Range('D50').Value = '0.06625'
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('D51').Value = '1.079'
$ws.Range('E51').Value = '  -4.69%  '
